$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Espinaca" (Vega Monumental
# Concepción). It belongs right before the current row 61, so insert a
# fresh row there and push the existing data (old rows 61-119) down by
# one (they end up as rows 62-120).
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new record's data.
$ws.Range("A61").Value = 11
$ws.Range("B61").Value = "Vega Monumental Concepción"
$ws.Range("C61").Value = "Bíobío"
$ws.Range("D61").Value = 45049
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = 100112012
$ws.Range("G61").Value = "Espinaca"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 90
$ws.Range("K61").Value = 8500
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = 9222
$ws.Range("N61").Value = "`$/cuna 10 kilos"
$ws.Range("O61").Value = "Región Metropolitana"
$ws.Range("P61").Value = 922
$ws.Range("Q61").Value = 10
$ws.Range("R61").Value = "Hortaliza"
